$d = $word.ActiveDocument

# First paragraph: "This is a Microsoft word document."
$para = $d.Paragraphs(1)
$r = $para.Range
$r.Text = "This is a Microsoft word document.  "

# Insert a new run after the existing text with the additional colored text.
$insertRange = $d.Range($para.Range.End - 1, $para.Range.End - 1)
$insertRange.Text = "(This is a change " + [char]0x2013 + " Version for branch alternate)"
$insertRange.Font.Color = 192
